# Re-populate the watchlist table (NSE ticker lists in columns B:F) with the
# new set of symbols, then drop the two now-unused trailing rows so the
# sheet's used range shrinks from A1:F20 down to A1:F18.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update changed cell values (rows 2-18)
$ws.Range("B2").Value = "NSE:ACC"
$ws.Range("C2").Value = "NSE:IFBIND"
$ws.Range("D2").Value = "NSE:CYIENT"
$ws.Range("E2").Value = "NSE:KOTAKBANK"
$ws.Range("F2").Value = "NSE:AUROPHARMA"
$ws.Range("B3").Value = "NSE:AUROPHARMA"
$ws.Range("C3").Value = "NSE:INDRAMEDCO"
$ws.Range("D3").Value = "NSE:JIOFIN"
$ws.Range("F3").Value = "NSE:BIOCON"
$ws.Range("B4").Value = "NSE:BAJAJELEC"
$ws.Range("C4").Value = "NSE:KAJARIACER"
$ws.Range("F4").Value = "NSE:HCLTECH"
$ws.Range("B5").Value = "NSE:BANSWRAS"
$ws.Range("C5").Value = "NSE:KALYANKJIL"
$ws.Range("F5").Value = "NSE:HDFCLIFE"
$ws.Range("B6").Value = "NSE:BIOCON"
$ws.Range("C6").Value = "NSE:MAXESTATES"
$ws.Range("F6").Value = "NSE:INDUSINDBK"
$ws.Range("B7").Value = "NSE:GIPCL"
$ws.Range("C7").Value = "NSE:RBLBANK"
$ws.Range("F7").Value = "NSE:JIOFIN"
$ws.Range("B8").Value = "NSE:GNFC"
$ws.Range("F8").Value = "NSE:LTIM"
$ws.Range("B9").Value = "NSE:HEXATRADEX"
$ws.Range("F9").Value = "NSE:MFSL"
$ws.Range("B10").Value = "NSE:INDIAMART"
$ws.Range("B11").Value = "NSE:INDUSINDBK"
$ws.Range("B12").Value = "NSE:INFY"
$ws.Range("B13").Value = "NSE:ITBEES"
$ws.Range("B14").Value = "NSE:JAMNAAUTO"
$ws.Range("B15").Value = "NSE:JSWENERGY"
$ws.Range("B16").Value = "NSE:MANORG"
$ws.Range("B17").Value = "NSE:MFSL"
$ws.Range("B18").Value = "NSE:PRECWIRE"

# Clear cells that no longer have a value
$ws.Range("E3").ClearContents()
$ws.Range("E4").ClearContents()
$ws.Range("E5").ClearContents()
$ws.Range("E6").ClearContents()
$ws.Range("E7").ClearContents()
$ws.Range("C8").ClearContents()
$ws.Range("E8").ClearContents()
$ws.Range("C9").ClearContents()
$ws.Range("E9").ClearContents()
$ws.Range("C10").ClearContents()
$ws.Range("E10").ClearContents()
$ws.Range("C11").ClearContents()
$ws.Range("E11").ClearContents()
$ws.Range("C12").ClearContents()
$ws.Range("E12").ClearContents()
$ws.Range("C13").ClearContents()
$ws.Range("E13").ClearContents()
$ws.Range("C14").ClearContents()
$ws.Range("E14").ClearContents()
$ws.Range("C15").ClearContents()
$ws.Range("C16").ClearContents()
$ws.Range("C17").ClearContents()
$ws.Range("C18").ClearContents()

# Remove the two trailing rows (19 and 20) so the used range shrinks to A1:F18
$ws.Rows.Item(19).Delete()
$ws.Rows.Item(19).Delete()
